# Bulk-booking-upload fixture refresh.
#
# The dates in the two sample booking rows (row 3 and row 4) were
# hard-coded and have drifted into the past relative to "today" for a
# Curve/Auslan e2e fixture, so they get bumped forward by ~10 years
# (the exact same time-of-day, just a later date) so the fixture keeps
# describing a bookable/future slot.
#
# D3/E3 (booking 1 start/end) : 2017-12-12 -> 2027-12-12
# D4/E4 (booking 2 start/end) : 2018-05-23 -> 2028-05-23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 46733.416666666664
$ws.Range("E3").Value = 46733.458333333336
$ws.Range("D4").Value = 46896.5
$ws.Range("E4").Value = 46896.583333333336

# Move the sheet's active selection onto the cells that were just
# edited (was previously parked on S3 with the view scrolled over to
# column N).
$ws.Range("D3:E4").Select()
